# Scheduled runner update: refresh currentAveragePrice / profit columns (H-N)
# for the leve rows below across all profession sheets, per market-board refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 2020
$ws.Range("J29").Value = 7998.3335
$ws.Range("L29").Value = 23995.0005
$ws.Range("N29").Value = -24557.0005
$ws.Range("H38").Value = 64.45
$ws.Range("I38").Value = 64.45
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 193.35
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H43").Value = 5885.857
$ws.Range("I43").Value = 6500.5
$ws.Range("J43").Value = 5640
$ws.Range("K43").Value = 6500.5
$ws.Range("L43").Value = 5640
$ws.Range("M43").Value = -6431.5
$ws.Range("N43").Value = -5778
$ws.Range("H58").Value = 1480.4706
$ws.Range("I58").Value = 270.5
$ws.Range("J58").Value = 2556
$ws.Range("K58").Value = 811.5
$ws.Range("L58").Value = 7668
$ws.Range("M58").Value = -661.5
$ws.Range("N58").Value = -7968
$ws.Range("H112").Value = 1923.375
$ws.Range("I112").Value = 930
$ws.Range("J112").Value = 2065.2856
$ws.Range("K112").Value = 2790
$ws.Range("L112").Value = 6195.8568
$ws.Range("M112").Value = -1682
$ws.Range("N112").Value = -8411.856800000001
$ws.Range("H137").Value = 4215.0625
$ws.Range("I137").Value = 2699.6667
$ws.Range("J137").Value = 4564.769
$ws.Range("K137").Value = 8099.000100000001
$ws.Range("L137").Value = 13694.307
$ws.Range("M137").Value = -5549.000100000001
$ws.Range("N137").Value = -18794.307

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1543.174
$ws.Range("I2").Value = 1333.4
$ws.Range("J2").Value = 1936.5
$ws.Range("K2").Value = 1333.4
$ws.Range("L2").Value = 1936.5
$ws.Range("M2").Value = -1220.4
$ws.Range("N2").Value = -2162.5
$ws.Range("H45").Value = 1560
$ws.Range("I45").Value = 1500
$ws.Range("J45").Value = 1600
$ws.Range("K45").Value = 1500
$ws.Range("L45").Value = 1600
$ws.Range("M45").Value = -1123
$ws.Range("N45").Value = -2354
$ws.Range("H61").Value = 2821.9092
$ws.Range("I61").Value = 2411.4375
$ws.Range("J61").Value = 3916.5
$ws.Range("K61").Value = 2411.4375
$ws.Range("L61").Value = 3916.5
$ws.Range("M61").Value = -2199.4375
$ws.Range("N61").Value = -4340.5
$ws.Range("H116").Value = 1543.174
$ws.Range("I116").Value = 1333.4
$ws.Range("J116").Value = 1936.5
$ws.Range("K116").Value = 1333.4
$ws.Range("L116").Value = 1936.5
$ws.Range("M116").Value = 960.5999999999999
$ws.Range("N116").Value = -6524.5
$ws.Range("H122").Value = 10586.479
$ws.Range("I122").Value = 10586.479
$ws.Range("K122").Value = 31759.437
$ws.Range("M122").Value = -29309.437
$ws.Range("H132").Value = 8903.5
$ws.Range("I132").Value = 12933.8
$ws.Range("J132").Value = 3865.625
$ws.Range("K132").Value = 38801.39999999999
$ws.Range("L132").Value = 11596.875
$ws.Range("M132").Value = -36271.39999999999
$ws.Range("N132").Value = -16656.875
$ws.Range("H134").Value = 49895
$ws.Range("J134").Value = 49895
$ws.Range("L134").Value = 49895
$ws.Range("N134").Value = -60035
$ws.Range("H136").Value = 2821.9092
$ws.Range("I136").Value = 2411.4375
$ws.Range("J136").Value = 3916.5
$ws.Range("K136").Value = 7234.3125
$ws.Range("L136").Value = 11749.5
$ws.Range("M136").Value = -4684.3125
$ws.Range("N136").Value = -16849.5
$ws.Range("H139").Value = 60932.082
$ws.Range("J139").Value = 60744.09
$ws.Range("L139").Value = 60744.09
$ws.Range("N139").Value = -71024.09

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1543.174
$ws.Range("I3").Value = 1333.4
$ws.Range("J3").Value = 1936.5
$ws.Range("K3").Value = 1333.4
$ws.Range("L3").Value = 1936.5
$ws.Range("M3").Value = -1219.4
$ws.Range("N3").Value = -2164.5
$ws.Range("H94").Value = 662.05554
$ws.Range("I94").Value = 613.4286
$ws.Range("J94").Value = 832.25
$ws.Range("K94").Value = 613.4286
$ws.Range("L94").Value = 832.25
$ws.Range("M94").Value = -162.4286
$ws.Range("N94").Value = -1734.25
$ws.Range("H134").Value = 2732.889
$ws.Range("I134").Value = 2227.56
$ws.Range("J134").Value = 3881.3635
$ws.Range("K134").Value = 6682.68
$ws.Range("L134").Value = 11644.0905
$ws.Range("M134").Value = -4147.68
$ws.Range("N134").Value = -16714.0905

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4716.5454
$ws.Range("I31").Value = 2283.3333
$ws.Range("K31").Value = 2283.3333
$ws.Range("M31").Value = -1988.3333
$ws.Range("H34").Value = 4716.5454
$ws.Range("I34").Value = 2283.3333
$ws.Range("K34").Value = 2283.3333
$ws.Range("M34").Value = -2081.3333
$ws.Range("H99").Value = 2390.4443
$ws.Range("I99").Value = 2342.8572
$ws.Range("J99").Value = 2557
$ws.Range("K99").Value = 2342.8572
$ws.Range("L99").Value = 2557
$ws.Range("M99").Value = -844.8571999999999
$ws.Range("N99").Value = -5553
$ws.Range("H126").Value = 2390.4443
$ws.Range("I126").Value = 2342.8572
$ws.Range("J126").Value = 2557
$ws.Range("K126").Value = 7028.571599999999
$ws.Range("L126").Value = 7671
$ws.Range("M126").Value = -4558.571599999999
$ws.Range("N126").Value = -12611
$ws.Range("H134").Value = 1652.4147
$ws.Range("I134").Value = 1338.826
$ws.Range("J134").Value = 2053.111
$ws.Range("K134").Value = 4016.478
$ws.Range("L134").Value = 6159.333
$ws.Range("M134").Value = -1481.478
$ws.Range("N134").Value = -11229.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 3963.0454
$ws.Range("J64").Value = 4376.3887
$ws.Range("L64").Value = 13129.1661
$ws.Range("N64").Value = -13669.1661
$ws.Range("H67").Value = 3963.0454
$ws.Range("J67").Value = 4376.3887
$ws.Range("L67").Value = 13129.1661
$ws.Range("N67").Value = -15001.1661
$ws.Range("H108").Value = 2821.6843
$ws.Range("I108").Value = 750.2857
$ws.Range("J108").Value = 4030
$ws.Range("K108").Value = 2250.8571
$ws.Range("L108").Value = 12090
$ws.Range("M108").Value = 629.1428999999998
$ws.Range("N108").Value = -17850
$ws.Range("H113").Value = 345514.62
$ws.Range("J113").Value = 946.7778
$ws.Range("L113").Value = 2840.3334
$ws.Range("N113").Value = -7180.3334
$ws.Range("H118").Value = 3332.6155
$ws.Range("J118").Value = 3332.6155
$ws.Range("L118").Value = 9997.8465
$ws.Range("N118").Value = -12483.8465
$ws.Range("H124").Value = 12440
$ws.Range("I124").Value = 1030
$ws.Range("J124").Value = 14070
$ws.Range("K124").Value = 3090
$ws.Range("L124").Value = 42210
$ws.Range("M124").Value = 1820
$ws.Range("N124").Value = -52030
$ws.Range("H125").Value = 3365
$ws.Range("I125").Value = 2000
$ws.Range("J125").Value = 3560
$ws.Range("K125").Value = 6000
$ws.Range("L125").Value = 10680
$ws.Range("M125").Value = -1080
$ws.Range("N125").Value = -20520
$ws.Range("H129").Value = 5000293
$ws.Range("I129").Value = 325.55554
$ws.Range("K129").Value = 976.66662
$ws.Range("M129").Value = 4023.33338
$ws.Range("H131").Value = 28664.879
$ws.Range("I131").Value = 694.4
$ws.Range("J131").Value = 40825.957
$ws.Range("K131").Value = 2083.2
$ws.Range("L131").Value = 122477.871
$ws.Range("M131").Value = 2956.8
$ws.Range("N131").Value = -132557.871
$ws.Range("H132").Value = 1791.4584
$ws.Range("I132").Value = 1422.9166
$ws.Range("J132").Value = 2160
$ws.Range("K132").Value = 12806.2494
$ws.Range("L132").Value = 19440
$ws.Range("M132").Value = -10276.2494
$ws.Range("N132").Value = -24500
$ws.Range("H133").Value = 4139.2
$ws.Range("I133").Value = 1970.2222
$ws.Range("J133").Value = 5359.25
$ws.Range("K133").Value = 5910.6666
$ws.Range("L133").Value = 16077.75
$ws.Range("M133").Value = -850.6665999999996
$ws.Range("N133").Value = -26197.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4167.7144
$ws.Range("I132").Value = 3708
$ws.Range("J132").Value = 4244.3335
$ws.Range("K132").Value = 11124
$ws.Range("L132").Value = 12733.0005
$ws.Range("M132").Value = -8594
$ws.Range("N132").Value = -17793.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2692.8
$ws.Range("I61").Value = 2799.2
$ws.Range("J61").Value = 2480
$ws.Range("K61").Value = 2799.2
$ws.Range("L61").Value = 2480
$ws.Range("M61").Value = -2597.2
$ws.Range("N61").Value = -2884
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H113").Value = 2692.8
$ws.Range("I113").Value = 2799.2
$ws.Range("J113").Value = 2480
$ws.Range("K113").Value = 2799.2
$ws.Range("L113").Value = 2480
$ws.Range("M113").Value = -629.1999999999998
$ws.Range("N113").Value = -6820
$ws.Range("H122").Value = 40910908
$ws.Range("I122").Value = 25001998
$ws.Range("K122").Value = 75005994
$ws.Range("M122").Value = -75003544
$ws.Range("H132").Value = 3633.4546
$ws.Range("I132").Value = 3629.2083
$ws.Range("K132").Value = 10887.6249
$ws.Range("M132").Value = -8357.624899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 51897.695
$ws.Range("I81").Value = 104104.09
$ws.Range("J81").Value = 4041.8333
$ws.Range("K81").Value = 208208.18
$ws.Range("L81").Value = 8083.6666
$ws.Range("M81").Value = -207147.18
$ws.Range("N81").Value = -10205.6666
$ws.Range("H84").Value = 51897.695
$ws.Range("I84").Value = 104104.09
$ws.Range("J84").Value = 4041.8333
$ws.Range("K84").Value = 1041040.9
$ws.Range("L84").Value = 40418.333
$ws.Range("M84").Value = -1035736.9
$ws.Range("N84").Value = -51026.333
$ws.Range("H122").Value = 13160631
$ws.Range("I122").Value = 13890971
$ws.Range("J122").Value = 14500
$ws.Range("K122").Value = 41672913
$ws.Range("L122").Value = 43500
$ws.Range("M122").Value = -41670463
$ws.Range("N122").Value = -48400
$ws.Range("H126").Value = 3500.5
$ws.Range("I126").Value = 3994.3076
$ws.Range("J126").Value = 2787.2222
$ws.Range("K126").Value = 11982.9228
$ws.Range("L126").Value = 8361.6666
$ws.Range("M126").Value = -9512.9228
$ws.Range("N126").Value = -13301.6666
$ws.Range("H132").Value = 5069.125
$ws.Range("I132").Value = 4911.1
$ws.Range("J132").Value = 5332.5
$ws.Range("K132").Value = 14733.3
$ws.Range("L132").Value = 15997.5
$ws.Range("M132").Value = -12203.3
$ws.Range("N132").Value = -21057.5
$ws.Range("H136").Value = 2123.524
$ws.Range("I136").Value = 1839.3103
$ws.Range("J136").Value = 2757.5386
$ws.Range("K136").Value = 5517.9309
$ws.Range("L136").Value = 8272.6158
$ws.Range("M136").Value = -2967.9309
$ws.Range("N136").Value = -13372.6158
